$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A87").Value = "Ulcer index"
$ws.Range("B87").Value = "Test ulcer index"
$ws.Range("C87").Value = "Ulcer_Index_test"

$ws.Range("C87").Select()
